$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2044
$ws.Range("I34").Value = 2044
$ws.Range("K34").Value = 2044
$ws.Range("M34").Value = -1841
$ws.Range("H36").Value = 2044
$ws.Range("I36").Value = 2044
$ws.Range("K36").Value = 2044
$ws.Range("M36").Value = -1329
$ws.Range("H116").Value = 4420.533
$ws.Range("I116").Value = 4761.8
$ws.Range("J116").Value = 4249.9
$ws.Range("K116").Value = 4761.8
$ws.Range("L116").Value = 4249.9
$ws.Range("M116").Value = -1319.8
$ws.Range("N116").Value = -11133.9
$ws.Range("H125").Value = 24408.572
$ws.Range("I125").Value = 54647.168
$ws.Range("K125").Value = 491824.512
$ws.Range("M125").Value = -489364.512
$ws.Range("H132").Value = 2759.68
$ws.Range("I132").Value = 1528.9048
$ws.Range("K132").Value = 4586.7144
$ws.Range("M132").Value = -2056.7144
$ws.Range("H137").Value = 4301.3335
$ws.Range("I137").Value = 1542.1562
$ws.Range("J137").Value = 26374.75
$ws.Range("K137").Value = 4626.4686
$ws.Range("L137").Value = 79124.25
$ws.Range("M137").Value = -2076.4686
$ws.Range("N137").Value = -84224.25
$ws.Range("H138").Value = 2571.5483
$ws.Range("I138").Value = 1441
$ws.Range("J138").Value = 4136.923
$ws.Range("K138").Value = 4323
$ws.Range("L138").Value = 12410.769
$ws.Range("M138").Value = 817
$ws.Range("N138").Value = -22690.769
$ws.Range("H141").Value = 48339.05
$ws.Range("I141").Value = 53501.777
$ws.Range("K141").Value = 160505.331
$ws.Range("M141").Value = -155325.331
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32052.854
$ws.Range("I32").Value = 17152.365
$ws.Range("K32").Value = 17152.365
$ws.Range("M32").Value = -16865.365
$ws.Range("H110").Value = 2217.889
$ws.Range("I110").Value = 2259.0881
$ws.Range("K110").Value = 2259.0881
$ws.Range("M110").Value = -214.0880999999999
$ws.Range("H122").Value = 1820
$ws.Range("I122").Value = 1697.1765
$ws.Range("K122").Value = 5091.529500000001
$ws.Range("M122").Value = -2641.529500000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1043.875
$ws.Range("I94").Value = 891.4286
$ws.Range("K94").Value = 891.4286
$ws.Range("M94").Value = -440.4286
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2258.4546
$ws.Range("I31").Value = 1829
$ws.Range("K31").Value = 1829
$ws.Range("M31").Value = -1534
$ws.Range("H34").Value = 2258.4546
$ws.Range("I34").Value = 1829
$ws.Range("K34").Value = 1829
$ws.Range("M34").Value = -1627
$ws.Range("H105").Value = 1596.8334
$ws.Range("J105").Value = 1619.5
$ws.Range("L105").Value = 1619.5
$ws.Range("N105").Value = -5113.5
$ws.Range("H134").Value = 2112.5652
$ws.Range("I134").Value = 2241.5264
$ws.Range("K134").Value = 6724.5792
$ws.Range("M134").Value = -4189.5792
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 612.4
$ws.Range("I5").Value = 431.33334
$ws.Range("J5").Value = 884
$ws.Range("K5").Value = 1294.00002
$ws.Range("L5").Value = 2652
$ws.Range("M5").Value = -1182.00002
$ws.Range("N5").Value = -2876
$ws.Range("H21").Value = 48.333332
$ws.Range("J21").Value = 48.333332
$ws.Range("L21").Value = 144.999996
$ws.Range("N21").Value = -490.999996
$ws.Range("H34").Value = 1572.8462
$ws.Range("I34").Value = 198.66667
$ws.Range("J34").Value = 2750.7144
$ws.Range("K34").Value = 596.00001
$ws.Range("L34").Value = 8252.143199999999
$ws.Range("M34").Value = -512.00001
$ws.Range("N34").Value = -8420.143199999999
$ws.Range("H39").Value = 79942
$ws.Range("J39").Value = 3575
$ws.Range("L39").Value = 10725
$ws.Range("N39").Value = -11313
$ws.Range("H55").Value = 10502742
$ws.Range("I55").Value = 334133
$ws.Range("J55").Value = 13892278
$ws.Range("K55").Value = 1002399
$ws.Range("L55").Value = 41676834
$ws.Range("M55").Value = -1002222
$ws.Range("N55").Value = -41677188
$ws.Range("H68").Value = 3114.6667
$ws.Range("J68").Value = 3922
$ws.Range("L68").Value = 11766
$ws.Range("N68").Value = -13388
$ws.Range("H71").Value = 3114.6667
$ws.Range("J71").Value = 3922
$ws.Range("L71").Value = 35298
$ws.Range("N71").Value = -43410
$ws.Range("H107").Value = 842.0833
$ws.Range("I107").Value = 882.5
$ws.Range("K107").Value = 2647.5
$ws.Range("M107").Value = -727.5
$ws.Range("H132").Value = 1466
$ws.Range("I132").Value = 1251.8572
$ws.Range("J132").Value = 1765.8
$ws.Range("K132").Value = 11266.7148
$ws.Range("L132").Value = 15892.2
$ws.Range("M132").Value = -8736.7148
$ws.Range("N132").Value = -20952.2
$ws.Range("H133").Value = 8885.866
$ws.Range("I133").Value = 6829
$ws.Range("K133").Value = 20487
$ws.Range("M133").Value = -15427
$ws.Range("H134").Value = 5825
$ws.Range("I134").Value = 2488.6365
$ws.Range("K134").Value = 7465.9095
$ws.Range("M134").Value = -2395.9095
$ws.Range("H135").Value = 612.4
$ws.Range("I135").Value = 431.33334
$ws.Range("J135").Value = 884
$ws.Range("K135").Value = 3882.00006
$ws.Range("L135").Value = 7956
$ws.Range("M135").Value = -1347.00006
$ws.Range("N135").Value = -13026
$ws.Range("H136").Value = 5922.5386
$ws.Range("I136").Value = 3374.125
$ws.Range("K136").Value = 10122.375
$ws.Range("M136").Value = -5022.375
$ws.Range("H138").Value = 7677.3335
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 7677.3335
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 23032.0005
$ws.Range("M138").Value = $null
$ws.Range("N138").Value = -33312.00049999999
$ws.Range("H140").Value = 2725.9333
$ws.Range("I140").Value = 2206.3572
$ws.Range("K140").Value = 6619.071599999999
$ws.Range("M140").Value = -1439.071599999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 133449.84
$ws.Range("I122").Value = 148966.2
$ws.Range("J122").Value = 1560.75
$ws.Range("K122").Value = 446898.6
$ws.Range("L122").Value = 4682.25
$ws.Range("M122").Value = -444448.6
$ws.Range("N122").Value = -9582.25
$ws.Range("H132").Value = 2532.258
$ws.Range("I132").Value = 2428.36
$ws.Range("K132").Value = 7285.08
$ws.Range("M132").Value = -4755.08
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 847.2857
$ws.Range("I22").Value = 863.5
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 863.5
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -568.5
$ws.Range("N22").Value = -1340
$ws.Range("H27").Value = 847.2857
$ws.Range("I27").Value = 863.5
$ws.Range("J27").Value = 750
$ws.Range("K27").Value = 863.5
$ws.Range("L27").Value = 750
$ws.Range("M27").Value = -756.5
$ws.Range("N27").Value = -964
$ws.Range("H46").Value = 1774.4375
$ws.Range("I46").Value = 1864.3
$ws.Range("K46").Value = 1864.3
$ws.Range("M46").Value = -1676.3
$ws.Range("H47").Value = 19782.166
$ws.Range("I47").Value = 9064
$ws.Range("J47").Value = 21925.8
$ws.Range("K47").Value = 9064
$ws.Range("L47").Value = 21925.8
$ws.Range("M47").Value = -8574
$ws.Range("N47").Value = -22905.8
$ws.Range("H52").Value = 19782.166
$ws.Range("I52").Value = 9064
$ws.Range("J52").Value = 21925.8
$ws.Range("K52").Value = 9064
$ws.Range("L52").Value = 21925.8
$ws.Range("M52").Value = -8831
$ws.Range("N52").Value = -22391.8
$ws.Range("H132").Value = 558519
$ws.Range("I132").Value = 791406.3
$ws.Range("K132").Value = 2374218.9
$ws.Range("M132").Value = -2371688.9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 25123.75
$ws.Range("J39").Value = 28498.334
$ws.Range("L39").Value = 28498.334
$ws.Range("N39").Value = -29324.334
$ws.Range("H43").Value = 9999
$ws.Range("I43").Value = 9999
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 9999
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -9850
$ws.Range("N43").Value = $null
$ws.Range("H51").Value = 32747.5
$ws.Range("J51").Value = 32747.5
$ws.Range("L51").Value = 32747.5
$ws.Range("N51").Value = -33767.5
$ws.Range("H52").Value = 18686.75
$ws.Range("I52").Value = 9999
$ws.Range("J52").Value = 19927.857
$ws.Range("K52").Value = 9999
$ws.Range("L52").Value = 19927.857
$ws.Range("M52").Value = -9773
$ws.Range("N52").Value = -20379.857
$ws.Range("H122").Value = 4639.2173
$ws.Range("I122").Value = 4366.7617
$ws.Range("K122").Value = 13100.2851
$ws.Range("M122").Value = -10650.2851
$ws.Range("H132").Value = 27795.482
$ws.Range("I132").Value = 26484.111
$ws.Range("K132").Value = 79452.333
$ws.Range("M132").Value = -76922.333
